$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.373.73"
$ws.Range("E2").Value = "  -1.94%  "

$ws.Range("D3").Value = "2.890.17"
$ws.Range("E3").Value = "  -2.74%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "524.04"
$ws.Range("E5").Value = "  -2.97%  "

$ws.Range("D6").Value = "140.87"
$ws.Range("E6").Value = "  -6.59%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -3.88%  "

$ws.Range("D9").Value = "2.893.30"
$ws.Range("E9").Value = "  -2.86%  "

$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -5.96%  "

$ws.Range("D11").Value = "5.97"
$ws.Range("E11").Value = "  -2.45%  "

$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("D13").Value = "3.395.82"
$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("E14").Value = "  +2.07%  "

$ws.Range("D15").Value = "60.407.80"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").Value = "22.56"
$ws.Range("E16").Value = "  -4.61%  "

$ws.Range("D17").Value = "2.894.62"
$ws.Range("E17").Value = "  -2.60%  "

$ws.Range("D18").Value = "'0.0000140"
$ws.Range("E18").Value = "  -4.36%  "

$ws.Range("D19").Value = "4.97"
$ws.Range("E19").Value = "  -3.85%  "

$ws.Range("D20").Value = "'11.60"
$ws.Range("E20").Value = "  -3.33%  "

$ws.Range("D21").Value = "350.86"
$ws.Range("E21").Value = "  -7.95%  "

$ws.Range("D22").Value = "6.58"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "5.71"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").Value = "64.37"
$ws.Range("E25").Value = "  -1.74%  "

$ws.Range("D26").Value = "0.451"
$ws.Range("E26").Value = "  -3.98%  "

$ws.Range("D27").Value = "0.178"
$ws.Range("E27").Value = "  -5.99%  "

$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "'7.80"
$ws.Range("E29").Value = "  -5.24%  "

$ws.Range("D30").Value = "0.0₃0830"
$ws.Range("E30").Value = "  -11.83%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").Value = "1.67"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("D33").Value = "'19.50"
$ws.Range("E33").Value = "  -4.70%  "

$ws.Range("D34").Value = "'150.40"
$ws.Range("E34").Value = "  -6.55%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.29"
$ws.Range("E35").Value = "  -7.76%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "5.55"
$ws.Range("E36").Value = "  -6.16%  "

$ws.Range("D37").Value = "0.992"
$ws.Range("E37").Value = "  -7.56%  "

$ws.Range("E38").Value = "  -5.33%  "

$ws.Range("D39").Value = "37.56"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  -5.57%  "

$ws.Range("D41").Value = "3.69"
$ws.Range("E41").Value = "  -5.61%  "

$ws.Range("D42").Value = "2.280.62"
$ws.Range("E42").Value = "  -5.45%  "

$ws.Range("D43").Value = "0.644"
$ws.Range("E43").Value = "  -3.49%  "

$ws.Range("D44").Value = "0.0579"
$ws.Range("E44").Value = "  -2.03%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'20.20"
$ws.Range("E45").Value = "  -8.89%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "'4.90"
$ws.Range("E47").Value = "  -4.46%  "

$ws.Range("D48").Value = "0.0237"
$ws.Range("E48").Value = "  -4.00%  "

$ws.Range("E49").Value = "  -1.28%  "

$ws.Range("D50").Value = "0.0914"
$ws.Range("E50").Value = "  -4.06%  "

$ws.Range("D51").Value = "247.03"
$ws.Range("E51").Value = "  -7.51%  "
